$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 7.777315673816282
$ws.Range("D2").Value = 3.014987419875711
$ws.Range("E2").Value = 24.79174213415143
$ws.Range("F2").Value = 19.4802831355137
$ws.Range("G2").Value = 23.44602904748765
$ws.Range("H2").Value = 11.19472305127429
$ws.Range("M2").Value = 41.70069791426189
$ws.Range("B3").Value = 7.704587222163366
$ws.Range("D3").Value = 2.922400963657579
$ws.Range("E3").Value = 23.41756359132795
$ws.Range("F3").Value = 19.26235927293049
$ws.Range("G3").Value = 22.68261639027289
$ws.Range("H3").Value = 11.28055079420221
$ws.Range("M3").Value = 39.19863302170486
$ws.Range("B4").Value = 7.66143945997142
$ws.Range("D4").Value = 2.86348191332191
$ws.Range("E4").Value = 22.53752749670013
$ws.Range("F4").Value = 19.15247696827177
$ws.Range("G4").Value = 22.2443913885128
$ws.Range("H4").Value = 11.34274728644962
$ws.Range("M4").Value = 37.57545245196752
$ws.Range("B5").Value = 7.644254268980068
$ws.Range("D5").Value = 2.838948601769874
$ws.Range("E5").Value = 22.17016049896304
$ws.Range("F5").Value = 19.11368910538449
$ws.Range("G5").Value = 22.07393207464346
$ws.Range("H5").Value = 11.37041752415979
$ws.Range("M5").Value = 36.8922928694506
$ws.Range("B6").Value = 7.641425244824957
$ws.Range("D6").Value = 2.834843187600237
$ws.Range("E6").Value = 22.10864393996882
$ws.Range("F6").Value = 19.1076088383987
$ws.Range("G6").Value = 22.04613008023763
$ws.Range("H6").Value = 11.37515068608002
$ws.Range("M6").Value = 36.77754963350671
$ws.Range("B7").Value = 7.661206059281472
$ws.Range("D7").Value = 2.863153170619747
$ws.Range("E7").Value = 22.53260790701708
$ws.Range("F7").Value = 19.15192966475593
$ws.Range("G7").Value = 22.24205907525466
$ws.Range("H7").Value = 11.34311112866277
$ws.Range("M7").Value = 37.56632674596789
$ws.Range("B8").Value = 7.751937400829574
$ws.Range("D8").Value = 2.983491615962308
$ws.Range("E8").Value = 24.32566026181651
$ws.Range("F8").Value = 19.40017168395404
$ws.Range("G8").Value = 23.17675304581721
$ws.Range("H8").Value = 11.22230695496894
$ws.Range("M8").Value = 40.85609745338637
$ws.Range("B9").Value = 7.94097164859457
$ws.Range("D9").Value = 3.203193058163686
$ws.Range("E9").Value = 27.54090789973296
$ws.Range("F9").Value = 20.07637466502008
$ws.Range("G9").Value = 25.23037646582564
$ws.Range("H9").Value = 11.06359295284519
$ws.Range("M9").Value = 46.61364606264339
$ws.Range("B10").Value = 8.085451329462391
$ws.Range("D10").Value = 3.354924191119205
$ws.Range("E10").Value = 29.90284592238014
$ws.Range("F10").Value = 20.68599608206198
$ws.Range("G10").Value = 26.84416751075701
$ws.Range("H10").Value = 10.99843756063836
$ws.Range("M10").Value = 50.41810888943049
$ws.Range("B11").Value = 8.152137514425347
$ws.Range("D11").Value = 3.421909218974017
$ws.Range("E11").Value = 30.95042488299778
$ws.Range("F11").Value = 20.98672400805461
$ws.Range("G11").Value = 27.5950144144461
$ws.Range("H11").Value = 10.98075104685631
$ws.Range("M11").Value = 52.0564801407947
$ws.Range("B12").Value = 8.177506464247273
$ws.Range("D12").Value = 3.446986979336879
$ws.Range("E12").Value = 31.33843134555673
$ws.Range("F12").Value = 21.10384663933233
$ws.Range("G12").Value = 27.88127336141023
$ws.Range("H12").Value = 10.97583377974828
$ws.Range("M12").Value = 52.66364430201877
$ws.Range("B13").Value = 8.172037981163433
$ws.Range("D13").Value = 3.441598759264276
$ws.Range("E13").Value = 31.25525269206108
$ws.Range("F13").Value = 21.07848004805681
$ws.Range("G13").Value = 27.81954374781593
$ws.Range("H13").Value = 10.97681267606175
$ws.Range("M13").Value = 52.5334694828127
$ws.Range("B14").Value = 8.154222426954231
$ws.Range("D14").Value = 3.423978123067148
$ws.Range("E14").Value = 30.98252020248475
$ws.Range("F14").Value = 20.99629551132408
$ws.Range("G14").Value = 27.61852903615821
$ws.Range("H14").Value = 10.98031047977657
$ws.Range("M14").Value = 52.10669729938974
$ws.Range("B15").Value = 8.143324386557877
$ws.Range("D15").Value = 3.413147624285463
$ws.Range("E15").Value = 30.81433390120043
$ws.Range("F15").Value = 20.94637359736232
$ws.Range("G15").Value = 27.49563949602645
$ws.Range("H15").Value = 10.98268660098242
$ws.Range("M15").Value = 51.84356219988351
$ws.Range("B16").Value = 8.081110703749912
$ws.Range("D16").Value = 3.350505972393562
$ws.Range("E16").Value = 29.83316505286407
$ws.Range("F16").Value = 20.66680373636874
$ws.Range("G16").Value = 26.79539300850409
$ws.Range("H16").Value = 10.99983924927692
$ws.Range("M16").Value = 50.30918157756575
$ws.Range("B17").Value = 8.043175294606057
$ws.Range("D17").Value = 3.311558791719621
$ws.Range("E17").Value = 29.21571087519437
$ws.Range("F17").Value = 20.5012102529973
$ws.Range("G17").Value = 26.36974903867352
$ws.Range("H17").Value = 11.01346778993654
$ws.Range("M17").Value = 49.34425807527551
$ws.Range("B18").Value = 8.021447679120367
$ws.Range("D18").Value = 3.288965226064426
$ws.Range("E18").Value = 28.85483900868545
$ws.Range("F18").Value = 20.40817708398349
$ws.Range("G18").Value = 26.12653063886269
$ws.Range("H18").Value = 11.02242811239303
$ws.Range("M18").Value = 48.78056955745669
$ws.Range("B19").Value = 8.014107523198639
$ws.Range("D19").Value = 3.281282278316808
$ws.Range("E19").Value = 28.73166846064191
$ws.Range("F19").Value = 20.37706099311504
$ws.Range("G19").Value = 26.04447105418467
$ws.Range("H19").Value = 11.02565263721224
$ws.Range("M19").Value = 48.58822144923828
$ws.Range("B20").Value = 8.04720423898768
$ws.Range("D20").Value = 3.315724639745134
$ws.Range("E20").Value = 29.28203241347554
$ws.Range("F20").Value = 20.51860988256229
$ws.Range("G20").Value = 26.41489745084169
$ws.Range("H20").Value = 11.01190047963011
$ws.Range("M20").Value = 49.44787512055382
$ws.Range("B21").Value = 8.159452304100355
$ws.Range("D21").Value = 3.429161505237708
$ws.Range("E21").Value = 31.06286368148107
$ws.Range("F21").Value = 21.02034812966933
$ws.Range("G21").Value = 27.67752317523498
$ws.Range("H21").Value = 10.97923429898125
$ws.Range("M21").Value = 52.23241002358596
$ws.Range("B22").Value = 8.233480929621594
$ws.Range("D22").Value = 3.501622925059163
$ws.Range("E22").Value = 32.17612392017007
$ws.Range("F22").Value = 21.36710149353213
$ws.Range("G22").Value = 28.51377466185838
$ws.Range("H22").Value = 10.96828808908043
$ws.Range("M22").Value = 53.97504616500105
$ws.Range("B23").Value = 8.193916554058763
$ws.Range("D23").Value = 3.463100482961884
$ws.Range("E23").Value = 31.58656719909757
$ws.Range("F23").Value = 21.18035336612539
$ws.Range("G23").Value = 28.06658779514836
$ws.Range("H23").Value = 10.973158843452
$ws.Range("M23").Value = 53.0520226890192
$ws.Range("B24").Value = 8.045382496291325
$ws.Range("D24").Value = 3.313841890391434
$ws.Range("E24").Value = 29.2520667939351
$ws.Range("F24").Value = 20.5107367457955
$ws.Range("G24").Value = 26.39448117618743
$ws.Range("H24").Value = 11.01260556089415
$ws.Range("M24").Value = 49.40105770544187
$ws.Range("B25").Value = 7.888763509368227
$ws.Range("D25").Value = 3.14547768513315
$ws.Range("E25").Value = 26.70488086617026
$ws.Range("F25").Value = 19.87333447031228
$ws.Range("G25").Value = 24.65424196571466
$ws.Range("H25").Value = 11.09776374992619
$ws.Range("M25").Value = 45.13074669344798
